# Updated cryptos list with latest prices and volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.120.32'
$ws.Range("E2").Value = '  -3.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.644.21'
$ws.Range("E3").Value = '  -2.71%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '526.87'
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.98'
$ws.Range("E6").Value = '  -4.06%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  -1.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.63'
$ws.Range("E9").Value = '  -6.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.104'
$ws.Range("E10").Value = '  -2.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.337'
$ws.Range("E11").Value = '  -2.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.131'
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.107.71'
$ws.Range("E13").Value = '  -3.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '59.038.63'
$ws.Range("E14").Value = '  -3.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.05'
$ws.Range("E15").Value = '  -3.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000137'
$ws.Range("E16").Value = '  -1.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.633.23'
$ws.Range("E17").Value = '  -5.26%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '342.19'
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.47'
$ws.Range("E19").Value = '  -1.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.56'
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.37'
$ws.Range("E21").Value = '  -0.97%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.57'
$ws.Range("E23").Value = '  +3.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.422'
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.168'
$ws.Range("E25").Value = '  -2.94%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.996'
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.25'
$ws.Range("E27").Value = '  -1.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0802'
$ws.Range("E28").Value = '  -4.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.54'
$ws.Range("E29").Value = '  -3.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.61'
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.98'
$ws.Range("E32").Value = '  -1.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.39'
$ws.Range("E33").Value = '  -1.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.23'
$ws.Range("E34").Value = '  -1.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.20'
$ws.Range("E35").Value = '  -2.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.916'
$ws.Range("E36").Value = '  -0.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.878'
$ws.Range("E37").Value = '  -4.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.55'
$ws.Range("E38").Value = '  -2.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.47'
$ws.Range("E39").Value = '  -4.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.66'
$ws.Range("E40").Value = '  -1.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.609'
$ws.Range("E42").Value = '  -3.14%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '274.73'
$ws.Range("E43").Value = '  -3.29%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0974'
$ws.Range("E44").Value = '  -1.30%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.43'
$ws.Range("E45").Value = '  -4.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0541'
$ws.Range("E46").Value = '  -1.53%  '
$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.65'
$ws.Range("E47").Value = '  +1.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.053.67'
$ws.Range("E48").Value = '  -2.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0231'
$ws.Range("E49").Value = '  -0.98%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.19'
$ws.Range("E50").Value = '  -2.77%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.69'
$ws.Range("E51").Value = '  -6.35%  '
